$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44475
$ws.Range("K2").Value = 9000
$ws.Range("L2").Value = 10000
$ws.Range("M2").Value = 9500
$ws.Range("P2").Value = 3167
$ws.Range("D3").Value = 44419
$ws.Range("K3").Value = 10000
$ws.Range("M3").Value = 10000
$ws.Range("P3").Value = 3333
$ws.Range("D4").Value = 44391
$ws.Range("L4").Value = 10000
$ws.Range("M4").Value = 10000
$ws.Range("P4").Value = 3333
$ws.Range("D5").Value = 44433
$ws.Range("D6").Value = 44405
$ws.Range("L6").Value = 10500
$ws.Range("M6").Value = 10250
$ws.Range("P6").Value = 3417
$ws.Range("D7").Value = 44454
$ws.Range("J7").Value = 16
$ws.Range("K7").Value = 9500
$ws.Range("L7").Value = 10000
$ws.Range("M7").Value = 9750
$ws.Range("P7").Value = 3250
$ws.Range("D8").Value = 44279
$ws.Range("J8").Value = 16
$ws.Range("D9").Value = 44426
$ws.Range("D10").Value = 44356
$ws.Range("D11").Value = 44175
$ws.Range("J11").Value = 70
$ws.Range("K11").Value = 12000
$ws.Range("L11").Value = 12000
$ws.Range("M11").Value = 12000
$ws.Range("P11").Value = 4000
$ws.Range("D12").Value = 44468
$ws.Range("K12").Value = 10000
$ws.Range("L12").Value = 11000
$ws.Range("M12").Value = 10500
$ws.Range("P12").Value = 3500
$ws.Range("D13").Value = 44266
$ws.Range("J13").Value = 160
$ws.Range("L13").Value = 10000
$ws.Range("M13").Value = 10000
$ws.Range("P13").Value = 3333
$ws.Range("D14").Value = 44482
$ws.Range("K14").Value = 9000
$ws.Range("L14").Value = 10000
$ws.Range("M14").Value = 9500
$ws.Range("P14").Value = 3167
$ws.Range("D15").Value = 44524
$ws.Range("K15").Value = 9000
$ws.Range("M15").Value = 9500
$ws.Range("P15").Value = 3167
$ws.Range("D16").Value = 44349
$ws.Range("J16").Value = 12
$ws.Range("K16").Value = 10000
$ws.Range("L16").Value = 10000
$ws.Range("M16").Value = 10000
$ws.Range("P16").Value = 3333
$ws.Range("D17").Value = 44321
$ws.Range("J17").Value = 25
$ws.Range("K17").Value = 10000
$ws.Range("M17").Value = 10000
$ws.Range("P17").Value = 3333
$ws.Range("D18").Value = 44412
$ws.Range("J18").Value = 25
$ws.Range("L18").Value = 10500
$ws.Range("M18").Value = 10260
$ws.Range("P18").Value = 3420
$ws.Range("D19").Value = 44517
$ws.Range("K19").Value = 9000
$ws.Range("M19").Value = 9500
$ws.Range("P19").Value = 3167
$ws.Range("D21").Value = 44328
$ws.Range("J21").Value = 16
$ws.Range("D22").Value = 44195
$ws.Range("J22").Value = 30
$ws.Range("D23").Value = 44272
$ws.Range("J23").Value = 70
$ws.Range("K23").Value = 10000
$ws.Range("M23").Value = 10000
$ws.Range("P23").Value = 3333
$ws.Range("D24").Value = 44293
$ws.Range("D25").Value = 44510
$ws.Range("K25").Value = 9000
$ws.Range("M25").Value = 9500
$ws.Range("P25").Value = 3167
$ws.Range("D26").Value = 44363
$ws.Range("J26").Value = 16
$ws.Range("D27").Value = 44503
$ws.Range("K27").Value = 8000
$ws.Range("L27").Value = 9000
$ws.Range("M27").Value = 8500
$ws.Range("P27").Value = 2833
$ws.Range("D28").Value = 44461
$ws.Range("J28").Value = 16
$ws.Range("K28").Value = 9500
$ws.Range("M28").Value = 9750
$ws.Range("P28").Value = 3250
$ws.Range("D29").Value = 44384
$ws.Range("J29").Value = 25
$ws.Range("K29").Value = 10000
$ws.Range("L29").Value = 10500
$ws.Range("M29").Value = 10260
$ws.Range("P29").Value = 3420
$ws.Range("D30").Value = 44370
$ws.Range("D31").Value = 44435
$ws.Range("L31").Value = 10500
$ws.Range("M31").Value = 10250
$ws.Range("P31").Value = 3417
$ws.Range("D32").Value = 44489
$ws.Range("J32").Value = 16
$ws.Range("K32").Value = 9000
$ws.Range("L32").Value = 10000
$ws.Range("M32").Value = 9500
$ws.Range("P32").Value = 3167
$ws.Range("D33").Value = 44377
$ws.Range("J33").Value = 16
$ws.Range("L33").Value = 10500
$ws.Range("M33").Value = 10250
$ws.Range("P33").Value = 3417
$ws.Range("D34").Value = 44300
$ws.Range("L34").Value = 10000
$ws.Range("M34").Value = 10000
$ws.Range("P34").Value = 3333
$ws.Range("D35").Value = 44307
$ws.Range("J35").Value = 160
$ws.Range("K35").Value = 10000
$ws.Range("M35").Value = 10000
$ws.Range("P35").Value = 3333
$ws.Range("D36").Value = 44181
$ws.Range("J36").Value = 10
$ws.Range("L36").Value = 12000
$ws.Range("M36").Value = 11000
$ws.Range("P36").Value = 3667
$ws.Range("D37").Value = 44335
$ws.Range("K37").Value = 10000
$ws.Range("M37").Value = 10000
$ws.Range("P37").Value = 3333
$ws.Range("D38").Value = 44314
$ws.Range("D39").Value = 44342
$ws.Range("J39").Value = 17
$ws.Range("L39").Value = 10000
$ws.Range("M39").Value = 10000
$ws.Range("P39").Value = 3333
$ws.Range("D40").Value = 44447
$ws.Range("J40").Value = 16
$ws.Range("L40").Value = 10500
$ws.Range("M40").Value = 10250
$ws.Range("P40").Value = 3417
